# sprint_project.xlsx update:
#  - AUGMENTATION DATA (row4) finished -> add completion date, status DONE
#  - MODELLING WITH IMAGE AUGMENTATION (row6) finished -> add completion date, status DONE
#  - MODEL EVALUATION (row7) finished -> add completion date, status DONE
#  - New task "ADD IMAGE PROCESSING" inserted as row8 (ON GOING), pushing the
#    remaining rows (TEST MODEL PERFORMANCE ... DEPLOY) down by one
#  - Sheet1 becomes the active/selected sheet & tab, with F9 selected

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

# --- Mark finished tasks with their actual completion date + DONE status ---

# Row 4: AUGMENTATION DATA
$ws1.Range("D4").Value = "26 Oktober 2021"
$ws1.Range("F4").Value = "DONE"

# Row 6: MODELLING WITH IMAGE AUGMENTATION
$ws1.Range("D6").Value = "26 Oktober 2021"
$ws1.Range("F6").Value = "DONE"

# Row 7: MODEL EVALUATION
$ws1.Range("D7").Value = "8 Nopember 2021"
$ws1.Range("F7").Value = "DONE"

# --- Insert a new task row for "ADD IMAGE PROCESSING" above the current row 8 ---
$ws1.Rows.Item(8).Insert()
$ws1.Rows.Item(8).ClearFormats()

$ws1.Range("A8").Value = "ADD IMAGE PROCESSING "
$ws1.Range("B8").Value = "9 Nopember 2021"
$ws1.Range("C8").Value = "16 Nopember 2021"
$ws1.Range("E8").Value = "FADHLAN"
$ws1.Range("F8").Value = "ON GOING"

# --- Make Sheet1 the active tab/sheet, with F9 selected ---
$ws1.Activate() | Out-Null
$ws1.Range("F9").Select() | Out-Null
